$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from H1 into the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I/J data values per row (row -> [I, J])
$data = @{
    2 = @(11, 12)
    3 = @(9, 9)
    4 = @(5, 5)
    5 = @(8, 8)
    6 = @(7, 7)
    7 = @(9, 9)
    8 = @(7, 8)
    9 = @(6, 6)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(10, 10)
    25 = @(9, 9)
    26 = @(10, 10)
    27 = @(9, 9)
    28 = @(7, 7)
    29 = @(8, 8)
    30 = @(9, 9)
    31 = @(8, 8)
    32 = @(8, 8)
    33 = @(9, 9)
    34 = @(8, 8)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(8, 9)
    38 = @(9, 9)
    39 = @(9, 9)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(9, 9)
    45 = @(7, 7)
    46 = @(8, 8)
    47 = @(9, 9)
    48 = @(8, 8)
    49 = @(10, 10)
    50 = @(9, 9)
    51 = @(6, 7)
    52 = @(8, 8)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(8, 8)
    57 = @(8, 8)
    58 = @(8, 8)
    59 = @(7, 8)
    60 = @(8, 9)
    61 = @(10, 10)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(6, 7)
    65 = @(8, 8)
    66 = @(6, 6)
    67 = @(9, 9)
    68 = @(8, 8)
    69 = @(5, 5)
    70 = @(8, 8)
    71 = @(6, 6)
    72 = @(6, 6)
    73 = @(5, 5)
    74 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
